$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 1
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 2
$ws.Range("F13").Value = -4
$ws.Range("F20").Value = 3
$ws.Range("F22").Value = 3
